$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of trade data (row 4), matching the formatting of row 3
$ws.Range("G3").Copy($ws.Range("G4"))

$ws.Range("A4").Value = 9983.91
$ws.Range("B4").Value = 9961
$ws.Range("C4").Value = 286.39
$ws.Range("D4").Value = 287.04000000000002
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.23
$ws.Range("G4").Value = 42608.639664351853
$ws.Range("H4").Value = $true
